# Applies the coinranking.com price/volume refresh described by the commit
# "Updated cryptos list on Thu Nov 23 10:13:35 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.584.76'
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").Value = '2.083.44'
$ws.Range("E3").Value = '  +3.65%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''236.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '''0.624'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("D7").Value = '''58.66'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.38%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.386'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.87%  '
$ws.Range("D10").Value = '''59.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.03%  '
$ws.Range("D11").Value = '''0.0765'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = '2.392.08'
$ws.Range("E13").Value = '  +3.71%  '
$ws.Range("D14").Value = '''14.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = '''21.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.07%  '
$ws.Range("E16").Value = '  +2.95%  '
$ws.Range("D17").Value = '''5.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '2.086.25'
$ws.Range("E18").Value = '  +3.68%  '
$ws.Range("D19").Value = '37.790.32'
$ws.Range("E19").Value = '  +3.01%  '
$ws.Range("D20").Value = '''6.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +15.67%  '
$ws.Range("D21").Value = '''69.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.70%  '
$ws.Range("D22").Value = '0.0₃0818'
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").Value = '''226.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  +3.70%  '
$ws.Range("D26").Value = '''2.45'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("D27").Value = '''167.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("E28").Value = '  +8.32%  '
$ws.Range("E29").Value = '  +4.16%  '
$ws.Range("D30").Value = '''0.131'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("D31").Value = '''19.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.66%  '
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").Value = '''4.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.01%  '
$ws.Range("D34").Value = '''0.0629'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.60%  '
$ws.Range("D35").Value = '''2.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.23%  '
$ws.Range("D36").Value = '''4.63'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.19%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '''5.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.83%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''3.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = '''4.56'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +21.04%  '
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("D43").Value = '''0.0960'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '1.473.61'
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''96.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.54%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '''1.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.11%  '
$ws.Range("E47").Value = '  +4.43%  '
$ws.Range("D48").Value = '''15.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.76%  '
$ws.Range("E49").Value = '  +4.14%  '
$ws.Range("D50").Value = '''7.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.97%  '
$ws.Range("E51").Value = '  +1.69%  '
